# Generate Report for Archive
#
# 1) The localization status moves on from "Ready for handoff" to
#    "In Translation" - update every cell that shows that status
#    (Overview!E2:F2, zh-cn!C2, de-de!C2).
# 2) The Status columns on the Overview / zh-cn / de-de sheets are
#    narrowed (report column re-sizing).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "Ready for handoff" -> "In Translation" ---
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# --- Narrower Status columns ---
# Target stored column width is ~13.41 characters; ColumnWidth is rounded
# to the nearest 1/6 character (pixel) internally, so 12.5 is the closest
# input that lands on the nearest achievable width.
$overview.Columns("E:F").ColumnWidth = 12.5
$zhcn.Columns("C:C").ColumnWidth = 12.5
$dede.Columns("C:C").ColumnWidth = 12.5
